$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Hunk 1: append ", getLoginStatus(http.Request, http.ResponseWriter)"
# as a NEW run (matching formatting) right after the existing
# "... getMods(http.Request, http.ResponseWriter)" run.
# -----------------------------------------------------------------

# Extend the existing run's text in place so the appended text
# inherits the exact run formatting (rFonts/color/sz/szCs) of the
# matched run.
$grow = $d.Content
$grow.Find.Execute("getMods(http.Request, http.ResponseWriter)", $true, $false, $false, $false, $false, $true, 1, $false, "getMods(http.Request, http.ResponseWriter), getLoginStatus(http.Request, http.ResponseWriter)", 1) | Out-Null

# Now split the appended tail off into its own run by toggling a
# character-formatting property on just that sub-range (round-trip
# back to the same value so the visible formatting is unchanged).
$tail = $d.Content
$tail.Find.Execute(", getLoginStatus(http.Request, http.ResponseWriter)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tail.Font.Bold = $true
$tail.Font.Bold = $false

# -----------------------------------------------------------------
# Hunk 2: merge the "Output" + ":" runs (immediately preceding the
# "If the controller is a Member controller, ..." paragraph) into a
# single "Output:" run.
# -----------------------------------------------------------------

$anchor = $d.Content
$anchor.Find.Execute("A client-server communication entity, communication.CommentReply, and a pointer to Server.") | Out-Null
$anchor.Collapse(0)
$anchor.Find.Execute("Output:", $true, $false, $false, $false, $false, $true, 1, $false, "Output:", 1) | Out-Null
